# Auto-generated Excel COM-interop script applying the leve market-data refresh
# described by the commit "chore: update Sheets via scheduled runner".
$wb = $excel.ActiveWorkbook

# ALC row 15
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1044.6538
$ws.Range("I15").Value = 1044.6538
$ws.Range("K15").Value = 3133.9614
$ws.Range("M15").Value = -2964.9614

# ALC row 32
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 7253
$ws.Range("I32").Value = 8499
$ws.Range("J32").Value = 6754.6
$ws.Range("K32").Value = 8499
$ws.Range("L32").Value = 6754.6
$ws.Range("M32").Value = -8173
$ws.Range("N32").Value = -7406.6

# ALC row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 37039584
$ws.Range("J86").Value = 3467.875
$ws.Range("L86").Value = 3467.875
$ws.Range("N86").Value = -5713.875

# ALC row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 37039584
$ws.Range("J89").Value = 3467.875
$ws.Range("L89").Value = 17339.375
$ws.Range("N89").Value = -28571.375

# ALC row 103
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 949
$ws.Range("J103").Value = 899
$ws.Range("L103").Value = 2697
$ws.Range("N103").Value = -3869

# ALC row 106
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 88237976
$ws.Range("I106").Value = 111113610
$ws.Range("K106").Value = 111113610
$ws.Range("M106").Value = -111112979

# ALC row 131
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 5976.6665
$ws.Range("I131").Value = 5072.778
$ws.Range("J131").Value = 11400
$ws.Range("K131").Value = 15218.334
$ws.Range("L131").Value = 34200
$ws.Range("M131").Value = -10178.334
$ws.Range("N131").Value = -44280

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1660.0834
$ws.Range("I132").Value = 1406.7407
$ws.Range("K132").Value = 4220.2221
$ws.Range("M132").Value = -1690.2221

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1928385.4
$ws.Range("I137").Value = 4285.515
$ws.Range("J137").Value = 5270243
$ws.Range("K137").Value = 12856.545
$ws.Range("L137").Value = 15810729
$ws.Range("M137").Value = -10306.545
$ws.Range("N137").Value = -15815829

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3923.7246
$ws.Range("J138").Value = 2973.0942
$ws.Range("L138").Value = 8919.2826
$ws.Range("N138").Value = -19199.2826

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 15385.125
$ws.Range("I141").Value = 2860.5
$ws.Range("K141").Value = 8581.5
$ws.Range("M141").Value = -3401.5

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 448559.6
$ws.Range("I74").Value = 1941.8485
$ws.Range("K74").Value = 1941.8485
$ws.Range("M74").Value = -1067.8485

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 448559.6
$ws.Range("I77").Value = 1941.8485
$ws.Range("K77").Value = 9709.2425
$ws.Range("M77").Value = -5341.2425

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1759.1017
$ws.Range("I132").Value = 1367.08
$ws.Range("K132").Value = 4101.24
$ws.Range("M132").Value = -1571.24

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 16365723
$ws.Range("I134").Value = 1832.0426
$ws.Range("K134").Value = 5496.1278
$ws.Range("M134").Value = -2961.1278

# CRP row 3
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 754113.3
$ws.Range("J3").Value = 5000
$ws.Range("L3").Value = 5000
$ws.Range("N3").Value = -5226

# CRP row 36
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 12500
$ws.Range("J36").Value = 15000
$ws.Range("L36").Value = 15000
$ws.Range("N36").Value = -15776

# CRP row 40
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H40").Value = 12500
$ws.Range("J40").Value = 15000
$ws.Range("L40").Value = 15000
$ws.Range("N40").Value = -15320

# CRP row 48
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5935.5
$ws.Range("I99").Value = 5756.8
$ws.Range("K99").Value = 5756.8
$ws.Range("M99").Value = -4258.8

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 5935.5
$ws.Range("I126").Value = 5756.8
$ws.Range("K126").Value = 17270.4
$ws.Range("M126").Value = -14800.4

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 15875146
$ws.Range("I132").Value = 1836.25
$ws.Range("K132").Value = 5508.75
$ws.Range("M132").Value = -2978.75

# CUL row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 26860826
$ws.Range("I4").Value = 52346510
$ws.Range("J4").Value = 100856.55
$ws.Range("K4").Value = 157039530
$ws.Range("L4").Value = 302569.65
$ws.Range("M4").Value = -157039418
$ws.Range("N4").Value = -302793.65

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 403.36667
$ws.Range("I5").Value = 408.8889
$ws.Range("J5").Value = 395.08334
$ws.Range("K5").Value = 1226.6667
$ws.Range("L5").Value = 1185.25002
$ws.Range("M5").Value = -1114.6667
$ws.Range("N5").Value = -1409.25002

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1495
$ws.Range("J68").Value = 1463.4375
$ws.Range("L68").Value = 4390.3125
$ws.Range("N68").Value = -6012.3125

# CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1495
$ws.Range("J71").Value = 1463.4375
$ws.Range("L71").Value = 13170.9375
$ws.Range("N71").Value = -21282.9375

# CUL row 129
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 28077856
$ws.Range("I129").Value = 1461.2858
$ws.Range("J129").Value = 44455750
$ws.Range("K129").Value = 4383.857400000001
$ws.Range("L129").Value = 133367250
$ws.Range("M129").Value = 616.1425999999992
$ws.Range("N129").Value = -133377250

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 6124709
$ws.Range("J131").Value = 5294729
$ws.Range("L131").Value = 15884187
$ws.Range("N131").Value = -15894267

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 403.36667
$ws.Range("I135").Value = 408.8889
$ws.Range("J135").Value = 395.08334
$ws.Range("K135").Value = 3680.0001
$ws.Range("L135").Value = 3555.75006
$ws.Range("M135").Value = -1145.0001
$ws.Range("N135").Value = -8625.75006

# CUL row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 3645.48
$ws.Range("J140").Value = 6142.5713
$ws.Range("L140").Value = 18427.7139
$ws.Range("N140").Value = -28787.7139

# GSM row 2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 74.5
$ws.Range("I2").Value = 58.53846
$ws.Range("J2").Value = 104.14286
$ws.Range("K2").Value = 58.53846
$ws.Range("L2").Value = 104.14286
$ws.Range("M2").Value = 54.46154
$ws.Range("N2").Value = -330.14286

# GSM row 40
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 5000
$ws.Range("J40").Value = 5000
$ws.Range("L40").Value = 5000
$ws.Range("N40").Value = -5302

# GSM row 64
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 60000
$ws.Range("J64").Value = 60000
$ws.Range("L64").Value = 60000
$ws.Range("N64").Value = -60496

# GSM row 67
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H67").Value = 60000
$ws.Range("J67").Value = 60000
$ws.Range("L67").Value = 60000
$ws.Range("N67").Value = -61716

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2979.5293
$ws.Range("J126").Value = 3972.3333
$ws.Range("L126").Value = 11916.9999
$ws.Range("N126").Value = -16856.9999

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 9497
$ws.Range("I22").Value = 7000
$ws.Range("K22").Value = 7000
$ws.Range("M22").Value = -6705

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 9497
$ws.Range("I27").Value = 7000
$ws.Range("K27").Value = 7000
$ws.Range("M27").Value = -6893

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 7044.6313
$ws.Range("J46").Value = 2300
$ws.Range("L46").Value = 2300
$ws.Range("N46").Value = -2676

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2774.5557
$ws.Range("I68").Value = 2717.2856
$ws.Range("J68").Value = 2975
$ws.Range("K68").Value = 2717.2856
$ws.Range("L68").Value = 2975
$ws.Range("M68").Value = -1968.2856
$ws.Range("N68").Value = -4473

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2774.5557
$ws.Range("I71").Value = 2717.2856
$ws.Range("J71").Value = 2975
$ws.Range("K71").Value = 13586.428
$ws.Range("L71").Value = 14875
$ws.Range("M71").Value = -9842.428
$ws.Range("N71").Value = -22363

# WVR row 13
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 1351.8334
$ws.Range("J13").Value = 1667
$ws.Range("L13").Value = 1667
$ws.Range("N13").Value = -1947

# WVR row 70
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 59998.75
$ws.Range("I70").Value = 59998
$ws.Range("K70").Value = 59998
$ws.Range("M70").Value = -59683

# WVR row 73
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H73").Value = 59998.75
$ws.Range("I73").Value = 59998
$ws.Range("K73").Value = 59998
$ws.Range("M73").Value = -58906

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1568.32
$ws.Range("I132").Value = 1195
$ws.Range("K132").Value = 3585
$ws.Range("M132").Value = -1055

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4335.727
$ws.Range("I136").Value = 4654.9546
$ws.Range("K136").Value = 13964.8638
$ws.Range("M136").Value = -11414.8638
